$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the comment attached to A2 ("cognate set" comment)
[void]$ws.Range("A2").Comment.Text("Cognate set comment")

# Row 2: swap B2/D2 string contents
#   B2 was "I" (bare pronoun gloss) -> becomes "WOMAN1"
#   D2 was "<em> (emphatic) {fict}" -> stays "<em> (emphatic) {fict}" (unchanged value, only shared-string index moves)
$ws.Range("B2").Value = "WOMAN1"
$ws.Range("D2").Value = "<em> (emphatic) {fict}"

# New row 3 with the second cognate-set entry
$ws.Range("A3").Value = "TG100"
$ws.Range("B3").Value = "PERSON1"
$ws.Range("D3").Value = "/am/ (description) {anysource}"

# Move the selection, matching the post-edit workbook's recorded cursor position
[void]$ws.Range("D14").Select()
